# Populate the newly-added "D" column labels for rows 16-21 (Med General /
# Med Cirjano / prot / bioq / orina / abdominal), reusing the same cell
# style already used by the other D-column labels (D6:D15), then leave the
# selection on E22 as in the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing label formatting (style index used by D6:D15) onto the
# new D16:D21 cells before filling in their text.
$ws.Range("D6").Copy()
$ws.Range("D16:D21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D16").Value = "Med General"
$ws.Range("D17").Value = "Med Cirjano"
$ws.Range("D18").Value = "prot"
$ws.Range("D19").Value = "bioq"
$ws.Range("D20").Value = "orina"
$ws.Range("D21").Value = "abdominal"

$ws.Range("E22").Select() | Out-Null
